$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.878.96"
$ws.Range("E2").Value = "  -0.93%  "
$ws.Range("D3").Value = "1.617.38"
$ws.Range("E3").Value = "  -1.55%  "
$ws.Range("E4").Value = "  -0.72%  "
$ws.Range("E5").Value = "  -1.23%  "
$ws.Range("E6").Value = "  -1.17%  "
$ws.Range("E7").Value = "  -0.71%  "
$ws.Range("E8").Value = "  -0.90%  "
$ws.Range("E9").Value = "  -1.38%  "
$ws.Range("D10").Value = "'18.31"
$ws.Range("E10").Value = "  -1.53%  "
$ws.Range("D11").Value = "'0.0790"
$ws.Range("E11").Value = "  -0.42%  "
$ws.Range("D12").Value = "1.840.89"
$ws.Range("E12").Value = "  -1.57%  "
$ws.Range("D13").Value = "1.612.42"
$ws.Range("E13").Value = "  -4.10%  "
$ws.Range("E14").Value = "  -1.61%  "
$ws.Range("E15").Value = "  -1.14%  "
$ws.Range("D16").Value = "25.883.96"
$ws.Range("E16").Value = "  -0.88%  "
$ws.Range("D17").Value = "'61.42"
$ws.Range("E17").Value = "  -1.53%  "
$ws.Range("E18").Value = "  -1.81%  "
$ws.Range("E19").Value = "  -0.72%  "
$ws.Range("D20").Value = "'191.11"
$ws.Range("E20").Value = "  +0.58%  "
$ws.Range("E21").Value = "  -0.95%  "
$ws.Range("D22").Value = "'9.50"
$ws.Range("E22").Value = "  -0.08%  "
$ws.Range("D23").Value = "'6.02"
$ws.Range("E23").Value = "  -1.72%  "
$ws.Range("E24").Value = "  +3.26%  "
$ws.Range("D25").Value = "'143.63"
$ws.Range("E25").Value = "  -0.31%  "
$ws.Range("E26").Value = "  -0.69%  "
$ws.Range("E27").Value = "  -3.02%  "
$ws.Range("E28").Value = "  -1.12%  "
$ws.Range("E29").Value = "  -0.11%  "
$ws.Range("D30").Value = "'1.23"
$ws.Range("E30").Value = "  -1.23%  "
$ws.Range("D31").Value = "'0.0477"
$ws.Range("E31").Value = "  -0.90%  "
$ws.Range("E32").Value = "  -1.46%  "
$ws.Range("D33").Value = "'3.10"
$ws.Range("E33").Value = "  -2.59%  "
$ws.Range("E34").Value = "  -1.80%  "
$ws.Range("D35").Value = "'1.49"
$ws.Range("E35").Value = "  -1.33%  "
$ws.Range("D36").Value = "1.125.92"
$ws.Range("E36").Value = "  +0.17%  "
$ws.Range("D37").Value = "'0.837"
$ws.Range("E37").Value = "  -4.48%  "
$ws.Range("E38").Value = "  -3.91%  "
$ws.Range("E39").Value = "  -1.47%  "
$ws.Range("E40").Value = "  -1.29%  "
$ws.Range("D41").Value = "'98.10"
$ws.Range("E41").Value = "  -0.81%  "
$ws.Range("D42").Value = "1.752.69"
$ws.Range("E42").Value = "  -1.29%  "
$ws.Range("E43").Value = "  -5.06%  "
$ws.Range("D44").Value = "'5.05"
$ws.Range("E44").Value = "  -4.28%  "
$ws.Range("E45").Value = "  -1.51%  "
$ws.Range("D46").Value = "'1.50"
$ws.Range("E46").Value = "  +1.95%  "
$ws.Range("D47").Value = "'53.99"
$ws.Range("E47").Value = "  -2.18%  "
$ws.Range("E48").Value = "  -0.62%  "
$ws.Range("E49").Value = "  -1.53%  "
$ws.Range("E50").Value = "  -0.76%  "
$ws.Range("D51").Value = "'7.44"
$ws.Range("E51").Value = "  -1.30%  "
